$word.Application.UserName = "Anders Meidahl"
$d = $word.ActiveDocument
$d.TrackRevisions = $true

# --- 1) Track-changes deletion of the two "noget " occurrences ---------
# "... hvor der indgår noget hardware og noget software. ..."
#                       ^^^^^^              ^^^^^^
# Locate both occurrences first (so positions stay valid), then delete
# starting from the later one so the earlier Range's offsets are not
# shifted by the second deletion.
$rng1 = $d.Content
$rng1.Find.Execute("noget ")

$rng2 = $d.Range($rng1.End, $d.Content.End)
$rng2.Find.Execute("noget ")

$rng2.Delete()
$rng1.Delete()

# Done recording tracked changes; turn the mode back off so it is not
# persisted as a document-wide setting (it wasn't part of the original
# authoring session state).
$d.TrackRevisions = $false

# --- 2) Move the "_GoBack" bookmark -------------------------------------
# Word re-seats its hidden "_GoBack" bookmark around the span of the most
# recent edit. Re-create it explicitly around "For en gennemgang af
# hvilke diagrammer der er benyttet, henvises der til dokumentationen for
# Pristjek220 ... hvor selve diagrammerne for projektet også kan findes."
# Adding a bookmark with the reserved name "_GoBack" replaces the
# existing one (Word only ever keeps a single "_GoBack" bookmark).
$bmStart = $d.Content
$bmStart.Find.Execute("For en gennemgang")

$bmEnd = $d.Content
$bmEnd.Find.Execute("hvor selve diagrammerne for projektet også kan findes.")

$bmRange = $d.Range($bmStart.Start, $bmEnd.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
